# Rename the three simulation worksheets
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Simulation 1").Name = "N"
$wb.Worksheets.Item("Simulation 2").Name = "2N"
$wb.Worksheets.Item("Simulation 3").Name = "3N"

# Update the selection on the active ("3N") sheet
$ws3 = $wb.Worksheets.Item("3N")
$ws3.Activate()
$ws3.Range("R49").Select()
